# Weekly update: add two new price entries for "Perejil" (Vega Modelo de Temuco)
# by inserting two new rows right before the old row 138. This pushes the
# existing rows 138:208 down to 140:210 (matching the new dimension A1:R210).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 138-139; everything below shifts down by 2.
$ws.Rows("138:139").Insert()

# Duplicate the row that is now at 140 (the old row 138) into the two new,
# still-empty rows so that all the "constant" columns (A,B,C,E,F,G,H,I,N,Q,R)
# are populated consistently with the rest of the table.
$ws.Range("A140:R140").Copy() | Out-Null
$ws.Range("A138:R138").PasteSpecial() | Out-Null
$ws.Range("A140:R140").Copy() | Out-Null
$ws.Range("A139:R139").PasteSpecial() | Out-Null
$excel.CutCopyMode = 0

# --- Row 138: new data point ---
$ws.Cells.Item(138, 4).Value = 44488   # D - Fecha
$ws.Cells.Item(138, 10).Value = 55     # J - Volumen
$ws.Cells.Item(138, 11).Value = 4000   # K - Precio minimo
$ws.Cells.Item(138, 12).Value = 4000   # L - Precio maximo
$ws.Cells.Item(138, 13).Value = 4000   # M - Precio promedio ponderado
$ws.Cells.Item(138, 15).Value = "Provincia de Cautín"  # O - Origen
$ws.Cells.Item(138, 16).Value = 1333   # P - Precio $/Kg

# --- Row 139: new data point ---
$ws.Cells.Item(139, 4).Value = 44488   # D - Fecha
$ws.Cells.Item(139, 10).Value = 115    # J - Volumen
$ws.Cells.Item(139, 11).Value = 3000   # K - Precio minimo
$ws.Cells.Item(139, 12).Value = 3000   # L - Precio maximo
$ws.Cells.Item(139, 13).Value = 3000   # M - Precio promedio ponderado
$ws.Cells.Item(139, 15).Value = "Región Metropolitana"  # O - Origen
$ws.Cells.Item(139, 16).Value = 1000   # P - Precio $/Kg
